$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Valor Mora" total (row 11) ---
$ws.Range("E11").Value = 160000

# --- Update "Cant. Periodos" count (row 13) ---
$ws.Range("F13").Value = 2

# --- Row 16: period 2507 -> 2506, Salario Basico 700000 -> 2000000 ---
$ws.Range("E16").Value = "2506"
$ws.Range("G16").Value = 2000000

# --- Row 17 becomes the new last data row: pick up row 18's (last-row) borders/format ---
$ws.Range("B18:J18").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New period for row 17 + updated Salario Basico value
$ws.Range("E17").Value = "2508"
$ws.Range("G17").Value = 2000000

# --- Remove the now-obsolete period row (old row 18, period 2505) ---
$ws.Rows(18).Delete()
